$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PLF 2021"

$ws.Range("A4").Value = "Solde public (en % du PIB)"
$ws.Range("B4").Value = "-3,0"
$ws.Range("C4").Value = "-10,2"
$ws.Range("D4").Value = "-6,7"

$ws.Range("A5").Value = "Croissance réelle du PIB (en %)"
$ws.Range("C5").Value = "-10,0"

$ws.Range("A6").Value = "Déflateur du PIB (en %)"

$ws.Range("A7").Value = "Commission européenne (prévisions de printemps 2020)"

$ws.Range("A8").Value = "Solde public (en % du PIB)"
$ws.Range("B8").Value = "-3,0"
$ws.Range("C8").Value = "-9,9"
$ws.Range("D8").Value = "-4,0"

$ws.Range("A9").Value = "Croissance réelle du PIB (en %)"
$ws.Range("C9").Value = "-8,2"

$ws.Range("A10").Value = "Déflateur du PIB (en %)"
